# Locate the shape that holds "Dobot was 81% accurate in all trails"
# (falls back to the known shape name if text search doesn't match).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "Dobot was 81% accurate in all trails") {
            $target = $sh
        }
    }
}
if ($target -eq $null) {
    $target = $s.Shapes.Item("TextBox 8")
}

$tr = $target.TextFrame.TextRange

# Original single run: "Dobot was 81% accurate in all trails"
# Becomes three runs (98% replaces 81%, with the middle run split out):
#   1: "Dobot "
#   2: "was 98% "   (was "was 81% ")
#   3: "accurate in all trails"
$mid = $tr.Characters(7, 8)
$mid.Text = "was 98% "

Write-Host "Updated text: $($tr.Text)"
